$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of
# silently converting to a Number and losing formatting (e.g. "1.00" -> 1).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated coin data (price + 1h volume change), and the two
# row swaps (15<->16, 41<->42, 50<->51) from the refreshed ranking feed.
$ws.Range("D2").Value = "27.164.02"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.627.08"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "214.77"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "20.40"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.634.57"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "4.15"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "0.546"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "64.81"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.142.68"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "217.22"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "0.999"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("D23").Value = "9.08"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "147.83"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "15.61"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "0.0508"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").Value = "1.354.44"
$ws.Range("E33").Value = "  +7.18%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").Value = "0.550"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "0.856"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.24"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "65.54"
$ws.Range("E42").Value = "  +6.13%  "
$ws.Range("D43").Value = "5.26"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "1.765.52"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "90.68"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "0.856"
$ws.Range("E47").Value = "  +29.00%  "
$ws.Range("D48").Value = "0.100"
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0968"
$ws.Range("E51").Value = "  -4.09%  "
